$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows 3-22 down to 4-23.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new data record.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "Femacal de La Calera"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = Get-Date -Year 2021 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Perfection"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 21000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 21521
$ws.Range("N3").Value = '$/saco 25 kilos'
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 861
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
